$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (142-144) ---
# Row 142: full record
$ws.Range("A142").Value = 43463
$ws.Range("B142").Value = 0.37152777777777773
$ws.Range("C142").Value = 0.44097222222222227

# Row 143: full record
$ws.Range("A143").Value = 43463
$ws.Range("B143").Value = 0.45624999999999999
$ws.Range("C143").Value = 0.55555555555555558

# Row 144: date + start time only (no end time)
$ws.Range("A144").Value = 43463
$ws.Range("B144").Value = 0.6479166666666667

# Rows 145-146 are formula-only (no A/B/C data)

# --- Extend the calculated columns (D/E/F) down through row 146 ---
# Re-filling the whole D139:F146 block as one fill turns it into a single
# shared-formula group spanning 139-146 (matching rows 139-141 picking up
# the new rows), exactly like Excel's own "fill down" behaviour.
$ws.Range("D139:D146").Formula = "=(C139-B139)* 1440"
$ws.Range("E139:E146").Formula = "=IF(C139>B139, (C139-B139)*1440, (B139-C139)*1440)"
$ws.Range("F139:F146").Formula = "=ABS((C139-B139)*1440)"

# --- Grow the table / autofilter to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F146"))

# --- Update the view: scroll position + active selection ---
$ws.Range("C144").Select() | Out-Null
